# Daily attendance processing - 2026-02-01 02:18:58
# Swap the order of recorders in the "Recorded By" column (G):
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
